# Auto-applied numeric refresh for Sheets (currentAveragePrice* / LevePrice* / LeveProfit* columns)
# Mirrors a scheduled market-data pull: only cached <v> numbers change, no formulas involved.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 7261
$ws.Cells.Item(76, 9).Value = 7396.75
$ws.Cells.Item(76, 11).Value = 7396.75
$ws.Cells.Item(76, 13).Value = -7081.75

$ws.Cells.Item(79, 8).Value = 7261
$ws.Cells.Item(79, 9).Value = 7396.75
$ws.Cells.Item(79, 11).Value = 7396.75
$ws.Cells.Item(79, 13).Value = -6304.75

$ws.Cells.Item(92, 8).Value = 66932.734
$ws.Cells.Item(92, 9).Value = 77194.16
$ws.Cells.Item(92, 11).Value = 77194.16
$ws.Cells.Item(92, 13).Value = -75946.16

$ws.Cells.Item(98, 8).Value = 3987.3333
$ws.Cells.Item(98, 9).Value = 3331.1667
$ws.Cells.Item(98, 10).Value = 5299.6665
$ws.Cells.Item(98, 11).Value = 3331.1667
$ws.Cells.Item(98, 12).Value = 5299.6665
$ws.Cells.Item(98, 13).Value = -1833.1667
$ws.Cells.Item(98, 14).Value = -8295.666499999999

$ws.Cells.Item(106, 8).Value = 2739.3333
$ws.Cells.Item(106, 9).Value = 2739.3333
$ws.Cells.Item(106, 11).Value = 2739.3333
$ws.Cells.Item(106, 13).Value = -2108.3333

$ws.Cells.Item(116, 8).Value = 378368.88
$ws.Cells.Item(116, 10).Value = 3849.4285
$ws.Cells.Item(116, 12).Value = 3849.4285
$ws.Cells.Item(116, 14).Value = -10733.4285

$ws.Cells.Item(122, 8).Value = 3987.3333
$ws.Cells.Item(122, 9).Value = 3331.1667
$ws.Cells.Item(122, 10).Value = 5299.6665
$ws.Cells.Item(122, 11).Value = 9993.500100000001
$ws.Cells.Item(122, 12).Value = 15898.9995
$ws.Cells.Item(122, 13).Value = -7543.500100000001
$ws.Cells.Item(122, 14).Value = -20798.9995

$ws.Cells.Item(125, 8).Value = 5854.7144
$ws.Cells.Item(125, 10).Value = 7806.8
$ws.Cells.Item(125, 12).Value = 70261.2
$ws.Cells.Item(125, 14).Value = -75181.2

$ws.Cells.Item(132, 8).Value = 6374.6665
$ws.Cells.Item(132, 9).Value = 3920.5
$ws.Cells.Item(132, 10).Value = 14964.25
$ws.Cells.Item(132, 11).Value = 11761.5
$ws.Cells.Item(132, 12).Value = 44892.75
$ws.Cells.Item(132, 13).Value = -9231.5
$ws.Cells.Item(132, 14).Value = -49952.75

$ws.Cells.Item(135, 8).Value = 37559.9
$ws.Cells.Item(135, 9).Value = 1323.7646
$ws.Cells.Item(135, 11).Value = 11913.8814
$ws.Cells.Item(135, 13).Value = -9378.8814

$ws.Cells.Item(138, 8).Value = 4007.4487
$ws.Cells.Item(138, 10).Value = 3516.606
$ws.Cells.Item(138, 12).Value = 10549.818
$ws.Cells.Item(138, 14).Value = -20829.818

$ws.Cells.Item(141, 8).Value = 15791.1
$ws.Cells.Item(141, 9).Value = 15791.1
$ws.Cells.Item(141, 11).Value = 47373.3
$ws.Cells.Item(141, 13).Value = -42193.3


$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(28, 8).Value = 30245.375
$ws.Cells.Item(28, 9).Value = 28423.285
$ws.Cells.Item(28, 11).Value = 28423.285
$ws.Cells.Item(28, 13).Value = -28231.285

$ws.Cells.Item(32, 8).Value = 6762.1084
$ws.Cells.Item(32, 9).Value = 4847.452
$ws.Cells.Item(32, 10).Value = 20739.1
$ws.Cells.Item(32, 11).Value = 4847.452
$ws.Cells.Item(32, 12).Value = 20739.1
$ws.Cells.Item(32, 13).Value = -4560.452
$ws.Cells.Item(32, 14).Value = -21313.1

$ws.Cells.Item(45, 8).Value = 773029.5600000001
$ws.Cells.Item(45, 9).Value = 1669074
$ws.Cells.Item(45, 11).Value = 1669074
$ws.Cells.Item(45, 13).Value = -1668697

$ws.Cells.Item(74, 8).Value = 128916.125
$ws.Cells.Item(74, 10).Value = 4233.3335
$ws.Cells.Item(74, 12).Value = 4233.3335
$ws.Cells.Item(74, 14).Value = -5981.3335

$ws.Cells.Item(77, 8).Value = 128916.125
$ws.Cells.Item(77, 10).Value = 4233.3335
$ws.Cells.Item(77, 12).Value = 21166.6675
$ws.Cells.Item(77, 14).Value = -29902.6675

$ws.Cells.Item(99, 8).Value = 30245.375
$ws.Cells.Item(99, 9).Value = 28423.285
$ws.Cells.Item(99, 11).Value = 28423.285
$ws.Cells.Item(99, 13).Value = -25428.285

$ws.Cells.Item(122, 8).Value = 5049.625
$ws.Cells.Item(122, 9).Value = 5500
$ws.Cells.Item(122, 11).Value = 16500
$ws.Cells.Item(122, 13).Value = -14050

$ws.Cells.Item(132, 8).Value = 257499.5
$ws.Cells.Item(132, 9).Value = 999999
$ws.Cells.Item(132, 11).Value = 2999997
$ws.Cells.Item(132, 13).Value = -2997467


$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3936.4443
$ws.Cells.Item(86, 9).Value = 3758.3076
$ws.Cells.Item(86, 10).Value = 4399.6
$ws.Cells.Item(86, 11).Value = 3758.3076
$ws.Cells.Item(86, 12).Value = 4399.6
$ws.Cells.Item(86, 13).Value = -2635.3076
$ws.Cells.Item(86, 14).Value = -6645.6

$ws.Cells.Item(89, 8).Value = 3936.4443
$ws.Cells.Item(89, 9).Value = 3758.3076
$ws.Cells.Item(89, 10).Value = 4399.6
$ws.Cells.Item(89, 11).Value = 18791.538
$ws.Cells.Item(89, 12).Value = 21998
$ws.Cells.Item(89, 13).Value = -13175.538
$ws.Cells.Item(89, 14).Value = -33230

$ws.Cells.Item(134, 8).Value = 2231.2727
$ws.Cells.Item(134, 9).Value = 1974.4
$ws.Cells.Item(134, 11).Value = 5923.200000000001
$ws.Cells.Item(134, 13).Value = -3388.200000000001


$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1829.6666
$ws.Cells.Item(16, 9).Value = 1829.6666
$ws.Cells.Item(16, 11).Value = 1829.6666
$ws.Cells.Item(16, 13).Value = -1542.6666

$ws.Cells.Item(113, 8).Value = 1829.6666
$ws.Cells.Item(113, 9).Value = 1829.6666
$ws.Cells.Item(113, 11).Value = 1829.6666
$ws.Cells.Item(113, 13).Value = 340.3334

$ws.Cells.Item(132, 8).Value = 2586.0625
$ws.Cells.Item(132, 9).Value = 2769
$ws.Cells.Item(132, 11).Value = 8307
$ws.Cells.Item(132, 13).Value = -5777

$ws.Cells.Item(134, 8).Value = 95762.91
$ws.Cells.Item(134, 9).Value = 116276.89
$ws.Cells.Item(134, 11).Value = 348830.67
$ws.Cells.Item(134, 13).Value = -346295.67


$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(10, 8).Value = 161.77777
$ws.Cells.Item(10, 9).Value = 161.77777
$ws.Cells.Item(10, 11).Value = 485.33331
$ws.Cells.Item(10, 13).Value = -346.33331

$ws.Cells.Item(23, 8).Value = 50000656
$ws.Cells.Item(23, 9).Value = 195
$ws.Cells.Item(23, 10).Value = 71429420
$ws.Cells.Item(23, 11).Value = 585
$ws.Cells.Item(23, 12).Value = 214288260
$ws.Cells.Item(23, 13).Value = -350
$ws.Cells.Item(23, 14).Value = -214288730

$ws.Cells.Item(29, 8).Value = 778.65216
$ws.Cells.Item(29, 10).Value = 789.9231
$ws.Cells.Item(29, 12).Value = 2369.7693
$ws.Cells.Item(29, 14).Value = -2923.7693

$ws.Cells.Item(81, 8).Value = 4999.3335
$ws.Cells.Item(81, 9).Value = 2999.5
$ws.Cells.Item(81, 10).Value = 8999
$ws.Cells.Item(81, 11).Value = 8998.5
$ws.Cells.Item(81, 12).Value = 26997
$ws.Cells.Item(81, 13).Value = -7875.5
$ws.Cells.Item(81, 14).Value = -29243

$ws.Cells.Item(84, 8).Value = 4999.3335
$ws.Cells.Item(84, 9).Value = 2999.5
$ws.Cells.Item(84, 10).Value = 8999
$ws.Cells.Item(84, 11).Value = 26995.5
$ws.Cells.Item(84, 12).Value = 80991
$ws.Cells.Item(84, 13).Value = -21379.5
$ws.Cells.Item(84, 14).Value = -92223

$ws.Cells.Item(122, 8).Value = 1151
$ws.Cells.Item(122, 10).Value = 1599.75
$ws.Cells.Item(122, 12).Value = 14397.75
$ws.Cells.Item(122, 14).Value = -19297.75

$ws.Cells.Item(132, 8).Value = 2156.125
$ws.Cells.Item(132, 9).Value = 1969.8
$ws.Cells.Item(132, 10).Value = 2466.6667
$ws.Cells.Item(132, 11).Value = 17728.2
$ws.Cells.Item(132, 12).Value = 22200.0003
$ws.Cells.Item(132, 13).Value = -15198.2
$ws.Cells.Item(132, 14).Value = -27260.0003

$ws.Cells.Item(137, 8).Value = 4292.222
$ws.Cells.Item(137, 10).Value = 6000
$ws.Cells.Item(137, 12).Value = 18000
$ws.Cells.Item(137, 14).Value = -28200

$ws.Cells.Item(140, 8).Value = 233594.16
$ws.Cells.Item(140, 9).Value = 252810.33
$ws.Cells.Item(140, 11).Value = 758430.99
$ws.Cells.Item(140, 13).Value = -753250.99

$ws.Cells.Item(141, 8).Value = 6170.6665
$ws.Cells.Item(141, 9).Value = 6170.6665
$ws.Cells.Item(141, 11).Value = 18511.9995
$ws.Cells.Item(141, 13).Value = -13331.9995


$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4921.0835
$ws.Cells.Item(70, 9).Value = 4358.125
$ws.Cells.Item(70, 10).Value = 6047
$ws.Cells.Item(70, 11).Value = 4358.125
$ws.Cells.Item(70, 12).Value = 6047
$ws.Cells.Item(70, 13).Value = -4088.125
$ws.Cells.Item(70, 14).Value = -6587

$ws.Cells.Item(73, 8).Value = 4921.0835
$ws.Cells.Item(73, 9).Value = 4358.125
$ws.Cells.Item(73, 10).Value = 6047
$ws.Cells.Item(73, 11).Value = 4358.125
$ws.Cells.Item(73, 12).Value = 6047
$ws.Cells.Item(73, 13).Value = -3422.125
$ws.Cells.Item(73, 14).Value = -7919

$ws.Cells.Item(107, 8).Value = 67591.53
$ws.Cells.Item(107, 10).Value = 1629.6
$ws.Cells.Item(107, 12).Value = 1629.6
$ws.Cells.Item(107, 14).Value = -5469.6


$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 8486.833000000001
$ws.Cells.Item(46, 9).Value = 1033.6666
$ws.Cells.Item(46, 10).Value = 9977.467000000001
$ws.Cells.Item(46, 11).Value = 1033.6666
$ws.Cells.Item(46, 12).Value = 9977.467000000001
$ws.Cells.Item(46, 13).Value = -845.6666
$ws.Cells.Item(46, 14).Value = -10353.467

$ws.Cells.Item(61, 8).Value = 3277.5
$ws.Cells.Item(61, 10).Value = 5873
$ws.Cells.Item(61, 12).Value = 5873
$ws.Cells.Item(61, 14).Value = -6277

$ws.Cells.Item(113, 8).Value = 3277.5
$ws.Cells.Item(113, 10).Value = 5873
$ws.Cells.Item(113, 12).Value = 5873
$ws.Cells.Item(113, 14).Value = -10213

$ws.Cells.Item(132, 8).Value = 57984.59
$ws.Cells.Item(132, 9).Value = 73403.766
$ws.Cells.Item(132, 10).Value = 5559.4
$ws.Cells.Item(132, 11).Value = 220211.298
$ws.Cells.Item(132, 12).Value = 16678.2
$ws.Cells.Item(132, 13).Value = -217681.298
$ws.Cells.Item(132, 14).Value = -21738.2


$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(69, 8).Value = 59442.57
$ws.Cells.Item(69, 10).Value = 59442.57
$ws.Cells.Item(69, 12).Value = 59442.57
$ws.Cells.Item(69, 14).Value = -60940.57

$ws.Cells.Item(72, 8).Value = 59442.57
$ws.Cells.Item(72, 10).Value = 59442.57
$ws.Cells.Item(72, 12).Value = 178327.71
$ws.Cells.Item(72, 14).Value = -185815.71

$ws.Cells.Item(132, 8).Value = 42983.773
$ws.Cells.Item(132, 9).Value = 52801.953
$ws.Cells.Item(132, 11).Value = 158405.859
$ws.Cells.Item(132, 13).Value = -155875.859

$ws.Cells.Item(136, 8).Value = 2611.709
$ws.Cells.Item(136, 9).Value = 2198.9185
$ws.Cells.Item(136, 10).Value = 5982.8335
$ws.Cells.Item(136, 11).Value = 6596.755500000001
$ws.Cells.Item(136, 12).Value = 17948.5005
$ws.Cells.Item(136, 13).Value = -4046.755500000001
$ws.Cells.Item(136, 14).Value = -23048.5005

